# Avance Migracion Sipro - progress updates across Controllers / Daos /
# Vistas / Avance sheets ("Cambios realizados en componente y carga de
# archivos project").

$wb = $excel.ActiveWorkbook

# --- Controllers ---------------------------------------------------------
$wsControllers = $wb.Worksheets.Item("Controllers")
$wsControllers.Activate()
$wsControllers.Range("D14").Select()
$wsControllers.Range("D14").Value = 0.05

# --- Daos ------------------------------------------------------------------
$wsDaos = $wb.Worksheets.Item("Daos")
$wsDaos.Activate()
$wsDaos.Range("C12").Value = 0.5
$wsDaos.Range("C13").Value = 0.5
$wsDaos.Range("C14").Value = 1
$wsDaos.Range("C15").Value = 1
$wsDaos.Range("C50").Value = 1
$wsDaos.Range("C94").Value = 1
$wsDaos.Range("C14").Select()

# --- Vistas ------------------------------------------------------------------
$wsVistas = $wb.Worksheets.Item("Vistas")
$wsVistas.Activate()
$wsVistas.Range("C8").Value = 0.1
$wsVistas.Range("C33").Value = 0.85
$wsVistas.Range("C37").Value = 1
$wsVistas.Range("C9").Select()

# --- Avance (resumen) -------------------------------------------------------
$wsAvance = $wb.Worksheets.Item("Avance")
$wsAvance.Activate()
$wsAvance.Range("D11").Value = 8.9
$wsAvance.Range("D12").Select()

Write-Output "Updated progress values across Controllers, Daos, Vistas and Avance"
